{"js": "// Replace each \"divided by\" expression in the worksheet table with its\n// updated problem, matching the canonical OOXML diff exactly (25 unique\n// `NNN\u00f7N=` cell texts swapped 1:1 for new ones).\nconst replacements = [\n  [\"674\u00f77=\", \"321\u00f78=\"],\n  [\"410\u00f74=\", \"207\u00f79=\"],\n  [\"745\u00f74=\", \"332\u00f78=\"],\n  [\"626\u00f78=\", \"733\u00f76=\"],\n  [\"453\u00f76=\", \"442\u00f77=\"],\n  [\"337\u00f78=\", \"195\u00f73=\"],\n  [\"217\u00f74=\", \"277\u00f75=\"],\n  [\"434\u00f72=\", \"668\u00f74=\"],\n  [\"220\u00f78=\", \"534\u00f78=\"],\n  [\"539\u00f73=\", \"826\u00f79=\"],\n  [\"254\u00f76=\", \"373\u00f78=\"],\n  [\"822\u00f76=\", \"378\u00f74=\"],\n  [\"400\u00f73=\", \"884\u00f72=\"],\n  [\"999\u00f73=\", \"322\u00f79=\"],\n  [\"598\u00f75=\", \"391\u00f74=\"],\n  [\"364\u00f73=\", \"771\u00f78=\"],\n  [\"593\u00f77=\", \"207\u00f76=\"],\n  [\"681\u00f72=\", \"244\u00f73=\"],\n  [\"159\u00f73=\", \"164\u00f73=\"],\n  [\"292\u00f74=\", \"706\u00f74=\"],\n  [\"246\u00f72=\", \"485\u00f78=\"],\n  [\"561\u00f73=\", \"242\u00f76=\"],\n  [\"427\u00f78=\", \"503\u00f74=\"],\n  [\"974\u00f74=\", \"504\u00f72=\"],\n  [\"248\u00f75=\", \"333\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"divided by\" expression in the worksheet table with its\n# updated problem, matching the canonical OOXML diff exactly (25 unique\n# `NNN\u00f7N=` cell texts swapped 1:1 for new ones).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"674\u00f77=\"; New=\"321\u00f78=\"}\n    @{Old=\"410\u00f74=\"; New=\"207\u00f79=\"}\n    @{Old=\"745\u00f74=\"; New=\"332\u00f78=\"}\n    @{Old=\"626\u00f78=\"; New=\"733\u00f76=\"}\n    @{Old=\"453\u00f76=\"; New=\"442\u00f77=\"}\n    @{Old=\"337\u00f78=\"; New=\"195\u00f73=\"}\n    @{Old=\"217\u00f74=\"; New=\"277\u00f75=\"}\n    @{Old=\"434\u00f72=\"; New=\"668\u00f74=\"}\n    @{Old=\"220\u00f78=\"; New=\"534\u00f78=\"}\n    @{Old=\"539\u00f73=\"; New=\"826\u00f79=\"}\n    @{Old=\"254\u00f76=\"; New=\"373\u00f78=\"}\n    @{Old=\"822\u00f76=\"; New=\"378\u00f74=\"}\n    @{Old=\"400\u00f73=\"; New=\"884\u00f72=\"}\n    @{Old=\"999\u00f73=\"; New=\"322\u00f79=\"}\n    @{Old=\"598\u00f75=\"; New=\"391\u00f74=\"}\n    @{Old=\"364\u00f73=\"; New=\"771\u00f78=\"}\n    @{Old=\"593\u00f77=\"; New=\"207\u00f76=\"}\n    @{Old=\"681\u00f72=\"; New=\"244\u00f73=\"}\n    @{Old=\"159\u00f73=\"; New=\"164\u00f73=\"}\n    @{Old=\"292\u00f74=\"; New=\"706\u00f74=\"}\n    @{Old=\"246\u00f72=\"; New=\"485\u00f78=\"}\n    @{Old=\"561\u00f73=\"; New=\"242\u00f76=\"}\n    @{Old=\"427\u00f78=\"; New=\"503\u00f74=\"}\n    @{Old=\"974\u00f74=\"; New=\"504\u00f72=\"}\n    @{Old=\"248\u00f75=\"; New=\"333\u00f74=\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$pair.Old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$pair.New, [ref]2) | Out-Null\n}\n"}
